$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.513.32"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.443.23"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'149.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.31%  "
$ws.Range("D7").Value = "3.444.26"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "'7.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "'0.392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "4.031.49"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'28.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.67%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "3.445.90"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "61.599.37"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'14.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "'9.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'388.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").Value = "'0.569"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").Value = "3.587.21"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").Value = "'72.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "'0.180"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").Value = "'7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "'1.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.03%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'24.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "'166.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").Value = "'26.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.78%  "
$ws.Range("D43").Value = "'0.793"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "'42.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "2.611.49"
$ws.Range("E48").Value = "  +5.68%  "
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("D50").Value = "'7.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.65%  "
$ws.Range("D51").Value = "'23.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
